# Update the PSSM numeric grid (B2:K21) on Sheet1 with the supplemental-figures values.
# (Sparse cells not covered by this diff keep the worksheet-wide "floor" value of -19.83129332503361
# from the prior revision; every data cell below is re-pointed to the new floor of -17.27846760795335
# or to its own updated figure.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -17.27846760795335
$ws.Range("C2").Value = -0.3802795714564407
$ws.Range("D2").Value = -17.27846760795335
$ws.Range("E2").Value = -17.27846760795335
$ws.Range("F2").Value = -17.27846760795335
$ws.Range("G2").Value = -17.27846760795335
$ws.Range("H2").Value = -17.27846760795335
$ws.Range("I2").Value = -17.27846760795335
$ws.Range("J2").Value = -17.27846760795335
$ws.Range("K2").Value = -17.27846760795335

# Row 3
$ws.Range("B3").Value = -17.27846760795335
$ws.Range("C3").Value = -17.27846760795335
$ws.Range("D3").Value = -17.27846760795335
$ws.Range("E3").Value = -17.27846760795335
$ws.Range("F3").Value = -17.27846760795335
$ws.Range("G3").Value = -17.27846760795335
$ws.Range("H3").Value = -17.27846760795335
$ws.Range("I3").Value = 1.378993903653756
$ws.Range("J3").Value = -17.27846760795335
$ws.Range("K3").Value = -17.27846760795335

# Row 4
$ws.Range("B4").Value = -17.27846760795335
$ws.Range("C4").Value = -0.4559341194688112
$ws.Range("D4").Value = 0.4172840332596945
$ws.Range("E4").Value = -17.27846760795335
$ws.Range("F4").Value = 3.90238912964075
$ws.Range("G4").Value = -17.27846760795335
$ws.Range("H4").Value = 1.045224912642486
$ws.Range("I4").Value = -17.27846760795335
$ws.Range("J4").Value = 2.085250094599819
$ws.Range("K4").Value = -17.27846760795335

# Row 5
$ws.Range("B5").Value = -17.27846760795335
$ws.Range("C5").Value = 0.2229108886976946
$ws.Range("D5").Value = -17.27846760795335
$ws.Range("E5").Value = -17.27846760795335
$ws.Range("F5").Value = -17.27846760795335
$ws.Range("G5").Value = 3.040045413023165
$ws.Range("H5").Value = -17.27846760795335
$ws.Range("I5").Value = -17.27846760795335
$ws.Range("J5").Value = -17.27846760795335
$ws.Range("K5").Value = -17.27846760795335

# Row 6
$ws.Range("B6").Value = -17.27846760795335
$ws.Range("C6").Value = -17.27846760795335
$ws.Range("D6").Value = -17.27846760795335
$ws.Range("E6").Value = -17.27846760795335
$ws.Range("F6").Value = -17.27846760795335
$ws.Range("G6").Value = -17.27846760795335
$ws.Range("H6").Value = -17.27846760795335
$ws.Range("I6").Value = -17.27846760795335
$ws.Range("J6").Value = -17.27846760795335
$ws.Range("K6").Value = -17.27846760795335

# Row 7
$ws.Range("B7").Value = 3.206580414555272
$ws.Range("C7").Value = -17.27846760795335
$ws.Range("D7").Value = -17.27846760795335
$ws.Range("E7").Value = -17.27846760795335
$ws.Range("F7").Value = -17.27846760795335
$ws.Range("G7").Value = -17.27846760795335
$ws.Range("H7").Value = -17.27846760795335
$ws.Range("I7").Value = -17.27846760795335
$ws.Range("J7").Value = -17.27846760795335
$ws.Range("K7").Value = -17.27846760795335

# Row 8
$ws.Range("B8").Value = -17.27846760795335
$ws.Range("C8").Value = -17.27846760795335
$ws.Range("D8").Value = -17.27846760795335
$ws.Range("E8").Value = 1.568372500039011
$ws.Range("F8").Value = -17.27846760795335
$ws.Range("G8").Value = -17.27846760795335
$ws.Range("H8").Value = -17.27846760795335
$ws.Range("I8").Value = -17.27846760795335
$ws.Range("J8").Value = -17.27846760795335
$ws.Range("K8").Value = -17.27846760795335

# Row 9
$ws.Range("B9").Value = 3.428717088613425
$ws.Range("C9").Value = -17.27846760795335
$ws.Range("D9").Value = -17.27846760795335
$ws.Range("E9").Value = -17.27846760795335
$ws.Range("F9").Value = -17.27846760795335
$ws.Range("G9").Value = -17.27846760795335
$ws.Range("H9").Value = -17.27846760795335
$ws.Range("I9").Value = -17.27846760795335
$ws.Range("J9").Value = -17.27846760795335
$ws.Range("K9").Value = -17.27846760795335

# Row 10
$ws.Range("B10").Value = -17.27846760795335
$ws.Range("C10").Value = -17.27846760795335
$ws.Range("D10").Value = -17.27846760795335
$ws.Range("E10").Value = -17.27846760795335
$ws.Range("F10").Value = -17.27846760795335
$ws.Range("G10").Value = -17.27846760795335
$ws.Range("H10").Value = -17.27846760795335
$ws.Range("I10").Value = 0.1547673529097875
$ws.Range("J10").Value = -17.27846760795335
$ws.Range("K10").Value = -17.27846760795335

# Row 11
$ws.Range("B11").Value = -17.27846760795335
$ws.Range("C11").Value = -17.27846760795335
$ws.Range("D11").Value = -17.27846760795335
$ws.Range("E11").Value = 2.154883509879888
$ws.Range("F11").Value = -17.27846760795335
$ws.Range("G11").Value = 1.490223954373763
$ws.Range("H11").Value = -17.27846760795335
$ws.Range("I11").Value = -17.27846760795335
$ws.Range("J11").Value = -17.27846760795335
$ws.Range("K11").Value = -17.27846760795335

# Row 12
$ws.Range("B12").Value = -17.27846760795335
$ws.Range("C12").Value = -17.27846760795335
$ws.Range("D12").Value = -17.27846760795335
$ws.Range("E12").Value = -17.27846760795335
$ws.Range("F12").Value = -17.27846760795335
$ws.Range("G12").Value = -17.27846760795335
$ws.Range("H12").Value = -17.27846760795335
$ws.Range("I12").Value = -17.27846760795335
$ws.Range("J12").Value = -17.27846760795335
$ws.Range("K12").Value = -17.27846760795335

# Row 13
$ws.Range("B13").Value = -17.27846760795335
$ws.Range("C13").Value = -17.27846760795335
$ws.Range("D13").Value = -17.27846760795335
$ws.Range("E13").Value = 2.110392657380081
$ws.Range("F13").Value = -17.27846760795335
$ws.Range("G13").Value = -17.27846760795335
$ws.Range("H13").Value = -17.27846760795335
$ws.Range("I13").Value = -17.27846760795335
$ws.Range("J13").Value = 1.387204874753335
$ws.Range("K13").Value = -17.27846760795335

# Row 14
$ws.Range("B14").Value = -17.27846760795335
$ws.Range("C14").Value = -17.27846760795335
$ws.Range("D14").Value = 1.016803592603806
$ws.Range("E14").Value = -17.27846760795335
$ws.Range("F14").Value = -17.27846760795335
$ws.Range("G14").Value = -17.27846760795335
$ws.Range("H14").Value = -17.27846760795335
$ws.Range("I14").Value = -17.27846760795335
$ws.Range("J14").Value = -17.27846760795335
$ws.Range("K14").Value = 4.321919473792721

# Row 15
$ws.Range("B15").Value = -17.27846760795335
$ws.Range("C15").Value = -17.27846760795335
$ws.Range("D15").Value = -0.3878589457010851
$ws.Range("E15").Value = -17.27846760795335
$ws.Range("F15").Value = -17.27846760795335
$ws.Range("G15").Value = -17.27846760795335
$ws.Range("H15").Value = -17.27846760795335
$ws.Range("I15").Value = -17.27846760795335
$ws.Range("J15").Value = -17.27846760795335
$ws.Range("K15").Value = -17.27846760795335

# Row 16
$ws.Range("B16").Value = -17.27846760795335
$ws.Range("C16").Value = -17.27846760795335
$ws.Range("D16").Value = -17.27846760795335
$ws.Range("E16").Value = -17.27846760795335
$ws.Range("F16").Value = -17.27846760795335
$ws.Range("G16").Value = -17.27846760795335
$ws.Range("H16").Value = -17.27846760795335
$ws.Range("I16").Value = -17.27846760795335
$ws.Range("J16").Value = 1.657828133441944
$ws.Range("K16").Value = -17.27846760795335

# Row 17
$ws.Range("B17").Value = -17.27846760795335
$ws.Range("C17").Value = 1.148636497840677
$ws.Range("D17").Value = -0.07178707166880911
$ws.Range("E17").Value = -17.27846760795335
$ws.Range("F17").Value = -17.27846760795335
$ws.Range("G17").Value = -17.27846760795335
$ws.Range("H17").Value = 2.161795309013912
$ws.Range("I17").Value = 1.753300891133579
$ws.Range("J17").Value = 2.577095919920291
$ws.Range("K17").Value = -17.27846760795335

# Row 18
$ws.Range("B18").Value = -17.27846760795335
$ws.Range("C18").Value = -17.27846760795335
$ws.Range("D18").Value = -17.27846760795335
$ws.Range("E18").Value = -17.27846760795335
$ws.Range("F18").Value = -17.27846760795335
$ws.Range("G18").Value = -17.27846760795335
$ws.Range("H18").Value = 1.995794508368834
$ws.Range("I18").Value = -0.07694669146371869
$ws.Range("J18").Value = 2.006449097165903
$ws.Range("K18").Value = -17.27846760795335

# Row 19
$ws.Range("B19").Value = -17.27846760795335
$ws.Range("C19").Value = -17.27846760795335
$ws.Range("D19").Value = 3.10586611491877
$ws.Range("E19").Value = -17.27846760795335
$ws.Range("F19").Value = -17.27846760795335
$ws.Range("G19").Value = -17.27846760795335
$ws.Range("H19").Value = 2.090288210489894
$ws.Range("I19").Value = 1.241755879384248
$ws.Range("J19").Value = -17.27846760795335
$ws.Range("K19").Value = -17.27846760795335

# Row 20
$ws.Range("B20").Value = -17.27846760795335
$ws.Range("C20").Value = 3.159419556307218
$ws.Range("D20").Value = 2.659062825996293
$ws.Range("E20").Value = -17.27846760795335
$ws.Range("F20").Value = 2.335316461939105
$ws.Range("G20").Value = -17.27846760795335
$ws.Range("H20").Value = 1.367574214518623
$ws.Range("I20").Value = 3.263272788942116
$ws.Range("J20").Value = -17.27846760795335
$ws.Range("K20").Value = -17.27846760795335

# Row 21
$ws.Range("B21").Value = -17.27846760795335
$ws.Range("C21").Value = 2.628483157397202
$ws.Range("D21").Value = -17.27846760795335
$ws.Range("E21").Value = 3.046628091874526
$ws.Range("F21").Value = -17.27846760795335
$ws.Range("G21").Value = 3.164367586157499
$ws.Range("H21").Value = 1.397506648380407
$ws.Range("I21").Value = -17.27846760795335
$ws.Range("J21").Value = -17.27846760795335
$ws.Range("K21").Value = -17.27846760795335
